$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.555.30"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.180.16"
$ws.Range("E3").Value = "  -4.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'572.24"
$ws.Range("D6").Value = "'169.68"
$ws.Range("E6").Value = "  -6.63%  "
$ws.Range("E7").Value = "  -6.14%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "3.190.63"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("D10").Value = "'0.120"
$ws.Range("E10").Value = "  -3.97%  "
$ws.Range("D11").Value = "'6.83"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("D13").Value = "3.743.07"
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "64.586.92"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "'25.39"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "3.194.30"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "'421.56"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'13.03"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "'7.18"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'70.36"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'5.68"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").Value = "'0.498"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -7.61%  "
$ws.Range("D29").Value = "'8.78"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.84"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.84"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'5.07"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").Value = "'157.31"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").Value = "2.720.58"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").Value = "  -4.64%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("E42").Value = "  -7.73%  "
$ws.Range("D43").Value = "'39.16"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "'0.718"
$ws.Range("E44").Value = "  -5.44%  "
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("D46").Value = "'5.56"
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").Value = "'292.41"
$ws.Range("E48").Value = "  -6.68%  "
$ws.Range("D49").Value = "'21.47"
$ws.Range("E50").Value = "  -5.56%  "

# Restore default (Normal) style for cells that needed a text-forcing apostrophe,
# so the quotePrefix formatting flag does not linger on these cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
